# Auto-generated edit applying scheduled-runner profit recalculations
# to the Hades_Profits workbook (ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1758.5797
$ws.Range("I15").Value = 1758.5797
$ws.Range("K15").Value = 5275.7391
$ws.Range("M15").Value = -5106.7391
$ws.Range("H43").Value = 1200.0571
$ws.Range("I43").Value = 400.2
$ws.Range("J43").Value = 1333.3667
$ws.Range("K43").Value = 400.2
$ws.Range("L43").Value = 1333.3667
$ws.Range("M43").Value = -331.2
$ws.Range("N43").Value = -1471.3667
$ws.Range("H125").Value = 799.8889
$ws.Range("I125").Value = 534.7143
$ws.Range("J125").Value = 1728
$ws.Range("K125").Value = 4812.428699999999
$ws.Range("L125").Value = 15552
$ws.Range("M125").Value = -2352.428699999999
$ws.Range("N125").Value = -20472
$ws.Range("H129").Value = 832
$ws.Range("J129").Value = 969.42426
$ws.Range("L129").Value = 2908.27278
$ws.Range("N129").Value = -12908.27278
$ws.Range("H135").Value = 52103.668
$ws.Range("I135").Value = 44479.78
$ws.Range("J135").Value = 63063
$ws.Range("K135").Value = 400318.02
$ws.Range("L135").Value = 567567
$ws.Range("M135").Value = -397783.02
$ws.Range("N135").Value = -572637
$ws.Range("H137").Value = 5265842.5
$ws.Range("I137").Value = 7694339
$ws.Range("J137").Value = 4100.5
$ws.Range("K137").Value = 23083017
$ws.Range("L137").Value = 12301.5
$ws.Range("M137").Value = -23080467
$ws.Range("N137").Value = -17401.5
$ws.Range("H138").Value = 2736584.8
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2736584.8
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 8209754.399999999
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -8220034.399999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5313501
$ws.Range("I32").Value = 6081965.5
$ws.Range("J32").Value = 30306.875
$ws.Range("K32").Value = 6081965.5
$ws.Range("L32").Value = 30306.875
$ws.Range("M32").Value = -6081678.5
$ws.Range("N32").Value = -30880.875
$ws.Range("H61").Value = 62626210
$ws.Range("I61").Value = 83417864
$ws.Range("K61").Value = 83417864
$ws.Range("M61").Value = -83417652
$ws.Range("H74").Value = 6811758.5
$ws.Range("I74").Value = 8367404
$ws.Range("J74").Value = 144707.14
$ws.Range("K74").Value = 8367404
$ws.Range("L74").Value = 144707.14
$ws.Range("M74").Value = -8366530
$ws.Range("N74").Value = -146455.14
$ws.Range("H77").Value = 6811758.5
$ws.Range("I77").Value = 8367404
$ws.Range("J77").Value = 144707.14
$ws.Range("K77").Value = 41837020
$ws.Range("L77").Value = 723535.7000000001
$ws.Range("M77").Value = -41832652
$ws.Range("N77").Value = -732271.7000000001
$ws.Range("H106").Value = 38612.727
$ws.Range("J106").Value = 38612.727
$ws.Range("L106").Value = 38612.727
$ws.Range("N106").Value = -41136.727
$ws.Range("H122").Value = 5850373
$ws.Range("I122").Value = 2612.2144
$ws.Range("J122").Value = 22224102
$ws.Range("K122").Value = 7836.6432
$ws.Range("L122").Value = 66672306
$ws.Range("M122").Value = -5386.6432
$ws.Range("N122").Value = -66677206
$ws.Range("H136").Value = 62626210
$ws.Range("I136").Value = 83417864
$ws.Range("K136").Value = 250253592
$ws.Range("M136").Value = -250251042

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 739.5599999999999
$ws.Range("I31").Value = 786.6842
$ws.Range("J31").Value = 728.50616
$ws.Range("K31").Value = 786.6842
$ws.Range("L31").Value = 728.50616
$ws.Range("M31").Value = -491.6842
$ws.Range("N31").Value = -1318.50616
$ws.Range("H34").Value = 739.5599999999999
$ws.Range("I34").Value = 786.6842
$ws.Range("J34").Value = 728.50616
$ws.Range("K34").Value = 786.6842
$ws.Range("L34").Value = 728.50616
$ws.Range("M34").Value = -584.6842
$ws.Range("N34").Value = -1132.50616
$ws.Range("H58").Value = 18547558
$ws.Range("I58").Value = 22668940
$ws.Range("J58").Value = 1346.2
$ws.Range("K58").Value = 22668940
$ws.Range("L58").Value = 1346.2
$ws.Range("M58").Value = -22668737
$ws.Range("N58").Value = -1752.2
$ws.Range("H136").Value = 18547558
$ws.Range("I136").Value = 22668940
$ws.Range("J136").Value = 1346.2
$ws.Range("K136").Value = 68006820
$ws.Range("L136").Value = 4038.6
$ws.Range("M136").Value = -68004270
$ws.Range("N136").Value = -9138.6
$ws.Range("H138").Value = 46722.223
$ws.Range("J138").Value = 46722.223
$ws.Range("L138").Value = 46722.223
$ws.Range("N138").Value = -57002.223

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 190000
$ws.Range("J37").Value = 190000
$ws.Range("L37").Value = 570000
$ws.Range("N37").Value = -570224
$ws.Range("H68").Value = 1242.4166
$ws.Range("I68").Value = 600.46155
$ws.Range("J68").Value = 2001.091
$ws.Range("K68").Value = 1801.38465
$ws.Range("L68").Value = 6003.272999999999
$ws.Range("M68").Value = -990.38465
$ws.Range("N68").Value = -7625.272999999999
$ws.Range("H71").Value = 1242.4166
$ws.Range("I71").Value = 600.46155
$ws.Range("J71").Value = 2001.091
$ws.Range("K71").Value = 5404.15395
$ws.Range("L71").Value = 18009.819
$ws.Range("M71").Value = -1348.15395
$ws.Range("N71").Value = -26121.819
$ws.Range("H76").Value = 3400
$ws.Range("J76").Value = 3466.6667
$ws.Range("L76").Value = 10400.0001
$ws.Range("N76").Value = -11166.0001
$ws.Range("H79").Value = 3400
$ws.Range("J79").Value = 3466.6667
$ws.Range("L79").Value = 10400.0001
$ws.Range("N79").Value = -13052.0001
$ws.Range("H107").Value = 818.4286
$ws.Range("I107").Value = 399.94
$ws.Range("J107").Value = 2428
$ws.Range("K107").Value = 1199.82
$ws.Range("L107").Value = 7284
$ws.Range("M107").Value = 720.1800000000001
$ws.Range("N107").Value = -11124
$ws.Range("H126").Value = 3067.1428
$ws.Range("I126").Value = 2892.5
$ws.Range("J126").Value = 3300
$ws.Range("K126").Value = 8677.5
$ws.Range("L126").Value = 9900
$ws.Range("M126").Value = -3737.5
$ws.Range("N126").Value = -19780
$ws.Range("H131").Value = 797.48486
$ws.Range("I131").Value = 470
$ws.Range("J131").Value = 920.2917
$ws.Range("K131").Value = 1410
$ws.Range("L131").Value = 2760.8751
$ws.Range("M131").Value = 3630
$ws.Range("N131").Value = -12840.8751
$ws.Range("H132").Value = 2032.72
$ws.Range("I132").Value = 1433.5834
$ws.Range("J132").Value = 2585.7693
$ws.Range("K132").Value = 12902.2506
$ws.Range("L132").Value = 23271.9237
$ws.Range("M132").Value = -10372.2506
$ws.Range("N132").Value = -28331.9237
$ws.Range("H141").Value = 21310
$ws.Range("I141").Value = 21310
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 63930
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -58750
$ws.Range("N141").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1290.4375
$ws.Range("I107").Value = 1062.1111
$ws.Range("J107").Value = 1584
$ws.Range("K107").Value = 1062.1111
$ws.Range("L107").Value = 1584
$ws.Range("M107").Value = 857.8888999999999
$ws.Range("N107").Value = -5424

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2582
$ws.Range("I7").Value = 2477.5
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 2477.5
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -2365.5
$ws.Range("N7").Value = -3224
$ws.Range("H40").Value = 7233.75
$ws.Range("I40").Value = 7410.2856
$ws.Range("J40").Value = 5998
$ws.Range("K40").Value = 7410.2856
$ws.Range("L40").Value = 5998
$ws.Range("M40").Value = -7274.2856
$ws.Range("N40").Value = -6270
$ws.Range("H126").Value = 2582
$ws.Range("I126").Value = 2477.5
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 7432.5
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -4962.5
$ws.Range("N126").Value = -13940

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1316.2222
$ws.Range("I126").Value = 870.8570999999999
$ws.Range("J126").Value = 2875
$ws.Range("K126").Value = 2612.5713
$ws.Range("L126").Value = 8625
$ws.Range("M126").Value = -142.5712999999996
$ws.Range("N126").Value = -13565
